$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 145
$ws.Range("F5").Value = 853
$ws.Range("F6").Value = 653
$ws.Range("F7").Value = 1220
$ws.Range("F9").Value = 803
$ws.Range("F10").Value = 689
$ws.Range("F11").Value = 262
$ws.Range("F14").Value = 722
$ws.Range("F15").Value = 920
$ws.Range("F16").Value = 9831
$ws.Range("F17").Value = 610
$ws.Range("F22").Value = 265
$ws.Range("F23").Value = 1755
$ws.Range("F24").Value = 26
$ws.Range("F26").Value = 482
$ws.Range("F27").Value = 179
$ws.Range("F29").Value = 266
$ws.Range("F32").Value = 64
$ws.Range("F33").Value = 96
$ws.Range("F35").Value = 176
$ws.Range("F36").Value = 192
$ws.Range("F37").Value = 168
$ws.Range("F38").Value = 38
$ws.Range("F25").Value = 291
$ws.Range("G25").Value = 238
# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 126
$ws.Range("F8").Value = 187
$ws.Range("F10").Value = 236
$ws.Range("F12").Value = 83
$ws.Range("F16").Value = 273
$ws.Range("G11").Value = 880
# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 813
# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 813
$ws.Range("F5").Value = 145
$ws.Range("F9").Value = 853
$ws.Range("F10").Value = 653
$ws.Range("F11").Value = 1220
$ws.Range("F13").Value = 126
$ws.Range("F14").Value = 803
$ws.Range("F15").Value = 689
$ws.Range("F16").Value = 262
$ws.Range("F19").Value = 920
$ws.Range("F20").Value = 9831
$ws.Range("F21").Value = 236
$ws.Range("F22").Value = 610
$ws.Range("F29").Value = 179
$ws.Range("F30").Value = 83
$ws.Range("F31").Value = 83
$ws.Range("F36").Value = 266
$ws.Range("F39").Value = 64
$ws.Range("F40").Value = 96
$ws.Range("F43").Value = 176
$ws.Range("F46").Value = 192
$ws.Range("F47").Value = 168
$ws.Range("B23").Value = "2024-08-17"
$ws.Range("C23").Value = "广州·第九届初物语动漫展内场—薄凉"
$ws.Range("D23").Value = "新港东路1000号 保利世贸博览馆"
$ws.Range("E23").Value = "2024.08.17 09:00-08.17 17:00"
$ws.Range("F23").Value = 43
$ws.Range("G23").Value = 138
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=89913"
$ws.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202407/gMZVWcQO1721983871813.jpeg"
$ws.Range("B24").Value = "2024-08-18"
$ws.Range("C24").Value = "广州·凹凸世界ONLY"
$ws.Range("D24").Value = "沙溪大道沙溪地铁站E2出口桥下 飞梦篮球公园(沙溪店)"
$ws.Range("E24").Value = "2024.08.18 10:00-08.18 17:00"
$ws.Range("F24").Value = 265
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=89715"
$ws.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202407/BnOHWZD91721638002542.jpeg"
$ws.Range("B25").Value = "2024-08-18"
$ws.Range("C25").Value = "广州·原神×崩坏×绝区零同人only"
$ws.Range("D25").Value = "西环路1号 广州岭南会展中心"
$ws.Range("E25").Value = "2024.08.18 10:00-08.18 17:00"
$ws.Range("F25").Value = 1755
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=87025"
$ws.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202405/lsOq4H701717169339283.png"
$ws.Range("B26").Value = "2024-08-18"
$ws.Range("C26").Value = "广州·宅舞联萌宅舞only同好会5.0·标记酱宅舞比赛3.0（免费活动）"
$ws.Range("D26").Value = "林和中路63号东方宝泰购物广场 东方宝泰购物广场B3层"
$ws.Range("E26").Value = "2024.08.18 12:00-08.18 19:15"
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 29.9
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=90059"
$ws.Range("I26").Value = "//i2.hdslb.com/bfs/openplatform/202407/eadPkAuz1722328902744.jpeg"
$ws.Range("B27").Value = "2024-08-23"
$ws.Range("C27").Value = "广州·LoveLiveOnly"
$ws.Range("D27").Value = "芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋"
$ws.Range("E27").Value = "2024.08.23 10:00-08.23 19:00"
$ws.Range("F27").Value = 482
$ws.Range("G27").Value = 68.8
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=87033"
$ws.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202406/a8shiH411717579829497.jpeg"
$ws.Range("B28").Value = "2024-08-24"
$ws.Range("C28").Value = "广州·LoveLive！电视动画播放十周年纪念巡演"
$ws.Range("D28").Value = "机场路1733号 久米空间LIVEHOUSE"
$ws.Range("E28").Value = "2024.08.24 12:30-08.25 18:30"
$ws.Range("F28").Value = 4439
$ws.Range("G28").Value = 880
$ws.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=86959"
$ws.Range("I28").Value = "//i1.hdslb.com/bfs/openplatform/202406/apzqBc5d1717661406596.jpeg"
